$p = $ppt.ActivePresentation

# Locate and remove the slide titled "Описание реализации" ("Описание
# реализации" describes the TicTacToe/Minesweeper/result/MainWindow
# classes) - it is dropped entirely from the deck.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
  $candidate = $p.Slides.Item($i)
  if ($candidate.Shapes.Item(1).TextFrame.TextRange.Text -eq "Описание реализации") {
    $candidate.Delete()
  }
}

# Find the slide titled "Вывод" and rework it into the "Результат" slide.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
  $s = $p.Slides.Item($i)
  if ($s.Shapes.Item(1).TextFrame.TextRange.Text -eq "Вывод") {
    $s.Shapes.Item(1).TextFrame.TextRange.Text = "Результат"
    $s.Shapes.Item(2).TextFrame.TextRange.Text = "Я выполнил задачу и доволен"
  }
}
